$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# Row 2 keeps a value in column C, but now references the new shared string.
$ws.Range("C2").Value = "PB-Controls And Appearance"

# Rows 3 through 63 in column C are cleared (style s="14" retained, content removed).
$ws.Range("C3:C63").Value = ""
